# Community amenities list (column A) gets refreshed with an updated,
# alphabetically sorted set of amenity names: several new amenities added,
# 'Golf Course' replaced by 'Grocery Service', and 'Storage Space' /
# 'Vintage Building' removed. Each amenity string literally begins and
# ends with a single quote character as part of its text content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    '''24 Hour Availability''',
    '''Adjacent to Open Space''',
    '''Adjacent to Walking / Biking Trails''',
    '''Basketball Court''',
    '''Bike Storage''',
    '''Boat Docks''',
    '''Breakfast/Coffee Concierge''',
    '''Business Center''',
    '''Cabana''',
    '''Car Charging Station''',
    '''Car Wash Area''',
    '''Clubhouse''',
    '''Clubhouse / Recreation Room''',
    '''Community Pool''',
    '''Community-Wide WiFi''',
    '''Composting''',
    '''Concierge''',
    '''Conference Room''',
    '''Controlled Access''',
    '''Corporate Suites''',
    '''Courtyard''',
    '''Cul de Sac''',
    '''Day Care''',
    '''Disposal Chutes''',
    '''Doorman''',
    '''Dry Cleaning Service''',
    '''Elevator''',
    '''Fenced Lot''',
    '''Fitness Center''',
    '''Full Scenic View''',
    '''Gameroom''',
    '''Gated Community''',
    '''Grocery Service''',
    '''Guest Apartment''',
    '''Health Club Discount''',
    '''Hearing Impaired Accessible''',
    '''Individual Leases Available''',
    '''Individual Locking Bedrooms''',
    '''Key Fob Entry''',
    '''Lake Access''',
    '''Laundry Service''',
    '''Lounge''',
    '''Maid Service''',
    '''Maintenance on site''',
    '''Meal Service''',
    '''Media Center/Movie Theatre''',
    '''Multi Use Room''',
    '''Near Parks''',
    '''Near Retail''',
    '''On-Site ATM''',
    '''On-Site Retail''',
    '''Online Services''',
    '''Package Service''',
    '''Partial Scenic View''',
    '''Pet Care''',
    '''Pet Play Area''',
    '''Pet Washing Station''',
    '''Picnic Area''',
    '''Planned Social Activities''',
    '''Playground''',
    '''Pond''',
    '''Private Bathroom''',
    '''Property Manager on Site''',
    '''Public Transportation''',
    '''Putting Greens''',
    '''Racquetball Court''',
    '''Recycling''',
    '''Renters Insurance Program''',
    '''Roof Terrace''',
    '''Roommate Matching''',
    '''Sauna''',
    '''Shuttle To Campus''',
    '''Shuttle to Train''',
    '''Spa''',
    '''Study Lounge''',
    '''Sundeck''',
    '''Tanning Salon''',
    '''Tennis Court''',
    '''Trash Pickup - Curbside''',
    '''Trash Pickup - Door to Door''',
    '''Vacuum System''',
    '''Video Patrol''',
    '''Vision Impaired Accessible''',
    '''Volleyball Court''',
    '''Walk To Campus''',
    '''Walking/Biking Trails''',
    '''Waterfront''',
    '''Waterfront View''',
    '''Wi-Fi at Pool and Clubhouse''',
    '''Zen Garden'''
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $cell = $ws.Cells.Item($i + 1, 1)
    # Setting .Value (or .Value2) directly to a string starting with "'"
    # makes Excel treat it as a "quote-prefix" text marker and strip it,
    # and it also tags the cell with a quotePrefix style. Instead, write
    # the literal text as a formula result, then paste-special as values
    # only, which preserves the leading apostrophe and keeps the cell
    # free of any extra number-format/style.
    $cell.Formula = '="' + $values[$i] + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163) # xlPasteValues
    $excel.CutCopyMode = 0
}

